$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Today" for this automated update run.
$today = Get-Date -Year 2026 -Month 2 -Day 10 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.UsedRange.Rows.Count
if (-not $lastRow -or $lastRow -lt 2) { $lastRow = 99 }

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $dVal = $dCell.Value2
    $fVal = $fCell.Value2

    if ($null -eq $dVal -or $null -eq $fVal) { continue }

    $fStr = [string]([int64]$fVal)
    if ($fStr.Length -ne 8) { continue }

    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)

    $startDate = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0
    $expireDate = $startDate.AddDays([int]$dVal)

    $remaining = [math]::Round($expireDate.ToOADate() - $today.ToOADate())

    if ($remaining -le 0) {
        # Supply ran out - restock today and reset the countdown.
        $eCell.Value2 = [int]$dVal
        $fCell.Value2 = [int]$today.ToString("yyyyMMdd")
    } else {
        $eCell.Value2 = $remaining
    }
}
